$d = $word.ActiveDocument

$replacements = @(
    @{old="284×7="; new="665×2="},
    @{old="551×4="; new="169×9="},
    @{old="378×4="; new="353×3="},
    @{old="168×2="; new="634×5="},
    @{old="313×7="; new="240×8="},
    @{old="976×3="; new="518×5="},
    @{old="895×9="; new="212×6="},
    @{old="877×9="; new="669×6="},
    @{old="582×7="; new="786×3="},
    @{old="662×8="; new="727×2="},
    @{old="412×7="; new="708×9="},
    @{old="716×7="; new="178×7="},
    @{old="239×7="; new="701×4="},
    @{old="842×8="; new="546×3="},
    @{old="313×5="; new="354×7="},
    @{old="589×3="; new="139×5="},
    @{old="125×9="; new="530×4="},
    @{old="442×9="; new="283×2="},
    @{old="466×5="; new="573×2="},
    @{old="869×4="; new="210×2="},
    @{old="648×4="; new="868×5="},
    @{old="616×3="; new="900×9="},
    @{old="134×6="; new="796×6="},
    @{old="983×6="; new="715×2="},
    @{old="360×9="; new="123×7="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
